$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5434
$ws.Range("L3").Value = 5890
$ws.Range("I4").Value = 1527
$ws.Range("L4").Value = 1443
$ws.Range("L5").Value = 350
$ws.Range("L6").Value = 4868
$ws.Range("I7").Value = 20968
$ws.Range("L7").Value = 17985

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L5").Value = 65
$ws.Range("L6").Value = 137
$ws.Range("L7").Value = 587
$ws.Range("L8").Value = 1186
$ws.Range("L9").Value = 104
$ws.Range("L10").Value = 119
$ws.Range("L11").Value = 295
$ws.Range("L16").Value = 37
$ws.Range("L20").Value = 443
$ws.Range("L25").Value = 109
$ws.Range("L29").Value = 1016
$ws.Range("L31").Value = 176
$ws.Range("L33").Value = 833
$ws.Range("L36").Value = 232
$ws.Range("L37").Value = 686
$ws.Range("L42").Value = 582
$ws.Range("L44").Value = 123
$ws.Range("L48").Value = 230
$ws.Range("L50").Value = 89
$ws.Range("L53").Value = 199
$ws.Range("L55").Value = 190
$ws.Range("L59").Value = 30
$ws.Range("I63").Value = 221
$ws.Range("L63").Value = 50
$ws.Range("L64").Value = 121
$ws.Range("L65").Value = 349
$ws.Range("L67").Value = 619
$ws.Range("L78").Value = 230
$ws.Range("L79").Value = 491
$ws.Range("L83").Value = 398
$ws.Range("L87").Value = 53
$ws.Range("L88").Value = 196
$ws.Range("L90").Value = 187
$ws.Range("L92").Value = 56
$ws.Range("L95").Value = 253
$ws.Range("I101").Value = 20968
$ws.Range("L101").Value = 17985

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 191
$ws.Range("L6").Value = 141
$ws.Range("L7").Value = 587

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 113
$ws.Range("L7").Value = 295

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L3").Value = 50
$ws.Range("L7").Value = 199

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 343
$ws.Range("L3").Value = 415
$ws.Range("L6").Value = 302
$ws.Range("L7").Value = 1186

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 126
$ws.Range("L7").Value = 398

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 229
$ws.Range("L3").Value = 293
$ws.Range("L7").Value = 833

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 80
$ws.Range("L7").Value = 253

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 206
$ws.Range("L3").Value = 241
$ws.Range("L7").Value = 686

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 109
$ws.Range("L5").Value = 6
$ws.Range("L7").Value = 349

# Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 72
$ws.Range("L7").Value = 176

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 183
$ws.Range("L3").Value = 237
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 619

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 303
$ws.Range("L3").Value = 390
$ws.Range("L7").Value = 1016

# Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 59
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 230

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 123

# Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 137

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 160
$ws.Range("L3").Value = 198
$ws.Range("L4").Value = 50
$ws.Range("L7").Value = 582

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 119

# Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 75
$ws.Range("L4").Value = 24
$ws.Range("L7").Value = 230

# Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L4").Value = 16
$ws.Range("L7").Value = 190

# Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 74
$ws.Range("L4").Value = 16

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 159
$ws.Range("L4").Value = 32
$ws.Range("L7").Value = 491

# Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 121

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 138
$ws.Range("L3").Value = 147
$ws.Range("L7").Value = 443

# Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 81
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 232

# East Side
$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L2").Value = 40
$ws.Range("L7").Value = 109

# Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 89

# Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 31
$ws.Range("L7").Value = 104

# Montclare
$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 30

# West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 56

# United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 196

# Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 65

# Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 187

# Ukrainian Village
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 53

# Bucktown
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 37
